$d = $word.ActiveDocument

function Set-ParaContentXml($paraIndex, [string]$innerXml) {
    $p = $d.Paragraphs.Item($paraIndex)
    $full = $p.Range
    # Sub-range excludes the trailing paragraph mark so pPr / the mark itself stay untouched.
    $sub = $d.Range($full.Start, $full.End - 1)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body><w:p>' + $innerXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $sub.InsertXML($pkg)
}

# --- Paragraph 1: "Minutes week " + "6" -> single run "Minutes week 6" ---
$p1Xml = '<w:r><w:rPr><w:rFonts w:asciiTheme="majorHAnsi" w:hAnsiTheme="majorHAnsi" w:cstheme="majorHAnsi"/><w:sz w:val="48"/><w:szCs w:val="48"/><w:lang w:val="en-US"/></w:rPr><w:t>Minutes week 6</w:t></w:r>'
Set-ParaContentXml 1 $p1Xml

# --- Paragraph 2: "Date and time:" run kept, " 1"+"9"+"-03-2018, 13:50" merged into one run ---
$p2Xml = '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>Date and time:</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> 19-03-2018, 13:50</w:t></w:r>'
Set-ParaContentXml 2 $p2Xml

# --- Paragraphs 4-6 (Chairman / Minute taker / Attendees) replaced together so the
#     proofErr markers that sit exactly on the paragraph-4 boundary are dropped too. ---
$chairmanXml = '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>Chairman</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>: Monika</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Kerulyte</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'

$minuteTakerXml = '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>Minute taker</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>: Ignas</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Kybransas</w:t></w:r>'

$attendeesXml = '<w:r><w:rPr><w:b/><w:lang w:val="en-US"/></w:rPr><w:t>Attendees</w:t></w:r>' +
    '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">: </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">Chung </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t>Kuah</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, Monika Kerulyte, Ignas Kybransas, </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:bCs/><w:color w:val="1D2129"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Teodor Genov</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:bCs/><w:color w:val="1D2129"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Vladimir Katrandjiev</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:eastAsia="lt-LT"/></w:rPr><w:t>,</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:bCs/><w:color w:val="1D2129"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Yoanna Borisova</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:cstheme="minorHAnsi"/><w:bCs/><w:color w:val="1D2129"/><w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/></w:rPr><w:t>Rostislav Tinchev</w:t></w:r>'

$p4 = $d.Paragraphs.Item(4)
$p6 = $d.Paragraphs.Item(6)
$block = $d.Range($p4.Range.Start, $p6.Range.End - 1)
$blockPkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + $chairmanXml + '</w:p>' +
    '<w:p><w:pPr><w:spacing w:after="120" w:line="257" w:lineRule="auto"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + $minuteTakerXml + '</w:p>' +
    '<w:p><w:pPr><w:spacing w:after="480" w:line="257" w:lineRule="auto"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' + $attendeesXml + '</w:p>' +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$block.InsertXML($blockPkg)

# --- Remove the stray bookmark that used to sit after "Mentor: thank you" ---
$d.Bookmarks.Item("_GoBack").Delete()
